# Insert a new weekly price record for "Terminal Hortofrutícola Agro Chillán -
# Frutilla" above the current row 158. This pushes the existing data rows
# (old 158..192) down to (159..193), matching the target sheet layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 158, shifting row 158 and everything below it down by one.
$ws.Rows("158:158").Insert()

# Populate the newly inserted row 158 with the new record's data.
$ws.Cells.Item(158, 1).Value = 7
$ws.Cells.Item(158, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(158, 3).Value = "Ñuble"
$ws.Cells.Item(158, 4).Value = 44543
$ws.Cells.Item(158, 5).Value = 16
$ws.Cells.Item(158, 6).Value = "Fruta"
$ws.Cells.Item(158, 7).Value = 100101
$ws.Cells.Item(158, 8).Value = "Berries"
$ws.Cells.Item(158, 9).Value = 100112025
$ws.Cells.Item(158, 10).Value = "Frutilla"
$ws.Cells.Item(158, 11).Value = "Sin especificar"
$ws.Cells.Item(158, 12).Value = "Primera"
$ws.Cells.Item(158, 13).Value = 160
$ws.Cells.Item(158, 14).Value = 7000
$ws.Cells.Item(158, 15).Value = 7500
$ws.Cells.Item(158, 16).Value = 7250
$ws.Cells.Item(158, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(158, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(158, 19).Value = 1036
$ws.Cells.Item(158, 20).Value = 7
